# Lit Review workbook update: add newly-read papers to the literature review
# table, push the old row-14/15 section boundary down and extend the sheet
# with several new entries + their usual "counter only" spacer rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Row 14 gains the thin/thick bottom-border look that marks the end of
#    a reference "block" (values/styles of the row itself are unchanged).
# ---------------------------------------------------------------------
$ws.Rows.Item(14).RowHeight = 15.75

# ---------------------------------------------------------------------
# 2. Row 15 becomes the first row of the next block, so its cells pick up
#    the "block-start" styling (same look already used by rows 20/22/26).
# ---------------------------------------------------------------------
$ws.Range("B20").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("C20").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D20").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E20").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("F20").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Rows.Item(15).RowHeight = 15.75

# ---------------------------------------------------------------------
# 3. Row 58 (last row of the previous block) gets the same bottom-border
#    treatment before the newly-appended material starts at row 59.
# ---------------------------------------------------------------------
$ws.Rows.Item(58).RowHeight = 15.75

# ---------------------------------------------------------------------
# 4. New reference #57 (row 59) - "block start" look, same as row 54.
# ---------------------------------------------------------------------
$ws.Range("B54").Copy()
$ws.Range("B59").PasteSpecial(-4122)
$ws.Range("C54").Copy()
$ws.Range("C59").PasteSpecial(-4122)
$ws.Range("D54").Copy()
$ws.Range("D59").PasteSpecial(-4122)
$ws.Range("E54").Copy()
$ws.Range("E59").PasteSpecial(-4122)
$ws.Range("F54").Copy()
$ws.Range("F59").PasteSpecial(-4122)

$ws.Range("B59").Value = 57
$ws.Range("C59").Value = "State of the Nation, Stroke Statistics February 2018"
$ws.Range("D59").Value = 5
$ws.Range("E59").Value = 1
$ws.Range("F59").Value = 2018
$ws.Rows.Item(59).RowHeight = 15.75

# ---------------------------------------------------------------------
# 5. New reference #58 (row 60).
# ---------------------------------------------------------------------
$ws.Range("C10").Copy()
$ws.Range("B60").PasteSpecial(-4122)
$ws.Range("C60").PasteSpecial(-4122)
$ws.Range("F60").PasteSpecial(-4122)
$ws.Range("D13").Copy()
$ws.Range("D60").PasteSpecial(-4122)
$ws.Range("E60").PasteSpecial(-4122)

$ws.Range("B60").Value = 58
$ws.Range("C60").Value = "Spasticity after stroke. Its occurrence and association with motor impairments and activity limitations."
$ws.Range("D60").Value = 4
$ws.Range("E60").Value = 2
$ws.Range("F60").Value = 2004
$ws.Rows.Item(60).RowHeight = 28.5

# ---------------------------------------------------------------------
# 6. Spacer / counter-only rows 61-64 (#59-#62).
# ---------------------------------------------------------------------
$ws.Range("B61").Value = 59
$ws.Range("B62").Value = 60
$ws.Range("B63").Value = 61
$ws.Range("B64").Value = 62

# ---------------------------------------------------------------------
# 7. New reference #63 (row 65).
# ---------------------------------------------------------------------
$ws.Range("C10").Copy()
$ws.Range("B65").PasteSpecial(-4122)
$ws.Range("C65").PasteSpecial(-4122)
$ws.Range("F65").PasteSpecial(-4122)
$ws.Range("D13").Copy()
$ws.Range("D65").PasteSpecial(-4122)
$ws.Range("E65").PasteSpecial(-4122)

$ws.Range("B65").Value = 63
$ws.Range("C65").Value = "Neurological Principles and Rehabilitation of Action Disorders: Rehabilitation interventions"
$ws.Range("D65").Value = 4
$ws.Range("E65").Value = 3
$ws.Range("F65").Value = 2011
$ws.Rows.Item(65).RowHeight = 30

# ---------------------------------------------------------------------
# 8. Spacer row 66 (#64).
# ---------------------------------------------------------------------
$ws.Range("B66").Value = 64

# ---------------------------------------------------------------------
# 9. New reference #65 (row 67).
# ---------------------------------------------------------------------
$ws.Range("C10").Copy()
$ws.Range("B67").PasteSpecial(-4122)
$ws.Range("C67").PasteSpecial(-4122)
$ws.Range("F67").PasteSpecial(-4122)
$ws.Range("D13").Copy()
$ws.Range("D67").PasteSpecial(-4122)
$ws.Range("E67").PasteSpecial(-4122)

$ws.Range("B67").Value = 65
$ws.Range("C67").Value = " The Science of Stroke: Mechanisms in search of Treatments"
$ws.Range("D67").Value = 3
$ws.Range("E67").Value = 4
$ws.Range("F67").Value = 2010

# ---------------------------------------------------------------------
# 10. Trailing counter-only rows 68-82 (#66-#80).
# ---------------------------------------------------------------------
$ws.Range("B68").Value = 66
$ws.Range("B69").Value = 67
$ws.Range("B70").Value = 68
$ws.Range("B71").Value = 69
$ws.Range("B72").Value = 70
$ws.Range("B73").Value = 71
$ws.Range("B74").Value = 72
$ws.Range("B75").Value = 73
$ws.Range("B76").Value = 74
$ws.Range("B77").Value = 75
$ws.Range("B78").Value = 76
$ws.Range("B79").Value = 77
$ws.Range("B80").Value = 78
$ws.Range("B81").Value = 79
$ws.Range("B82").Value = 80

# ---------------------------------------------------------------------
# 11. Update the view: scrolled down to the new material, selection on C70.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C70").Select()
